$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$none = "None"
$inputsValid = 'title = "Educated"' + "`n" + 'author = "Tara Westover"' + "`n" + 'genre = Genre.NON_FICTION'
$inputsBlankTitle = 'title = " "' + "`n" + 'author = "Tara Westover"' + "`n" + 'genre = Genre.NON_FICTION'
$inputsBlankAuthor = 'title = "Educated"' + "`n" + 'author = " "' + "`n" + 'genre = Genre.NON_FICTION'
$inputsBadGenre = 'title = "Educated"' + "`n" + 'author = "Tara Westover"' + "`n" + 'genre = "horror"'

# Developer name
$ws.Range("C3").Value = "Jashanpreet Sidhu "

# Preconditions for row 7
$ws.Range("E7").Value = $none

# Expected results
$ws.Range("G7").Value = "Attributes are set "
$ws.Range("G8").Value = "ValueError"
$ws.Range("G9").Value = "ValueError"
$ws.Range("G10").Value = "ValueError"

# Method Inputs
$ws.Range("F7").Value = $inputsValid
$ws.Range("F8").Value = $inputsBlankTitle
$ws.Range("F9").Value = $inputsBlankAuthor
$ws.Range("F10").Value = $inputsBadGenre

# Remaining Preconditions
$ws.Range("E8").Value = $none
$ws.Range("E9").Value = $none
$ws.Range("E10").Value = $none
$ws.Range("E11").Value = $inputsValid
$ws.Range("E12").Value = $inputsValid
$ws.Range("E13").Value = $inputsValid

# Remaining Method Inputs
$ws.Range("F11").Value = $none
$ws.Range("F12").Value = $none
$ws.Range("F13").Value = $none
